$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert "Date of Birth" column before PAN (col C)
$ws.Columns("C").Insert()

# Step 2: Insert "Correspondence Address" column before KYC Type (now col F)
$ws.Columns("F").Insert()

# Step 3: Insert "Bank Name" and "Branch Name" columns before Bank Account (now col I)
$ws.Columns("I:J").Insert()

# Step 4: Insert "Account Type" column after Bank Account Number (now col K) -> insert before L
$ws.Columns("L").Insert()

# Header row values for newly inserted columns
$ws.Range("C1").Value = "Date of Birth"
$ws.Range("F1").Value = "Correspondence Address"
$ws.Range("I1").Value = "Bank Name"
$ws.Range("J1").Value = "Branch Name"
$ws.Range("L1").Value = "Account Type"
$ws.Range("K1").Value = "Bank Account Number"

# Data row 2 (Investor 1)
$ws.Range("C2").Value = 27478
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").Value = "BOB"
$ws.Range("J2").Value = "Xyz"
$ws.Range("L2").Value = "Savings"

# Data row 3 (Investor 2)
$ws.Range("C3").Value = 15772
$ws.Range("I3").Value = "HDFC"
$ws.Range("J3").Value = "Abc"
$ws.Range("L3").Value = "Current"

# Reuse the same date style for C3 (copy format only from C2)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column width adjustments for the new columns
$ws.Columns("C").ColumnWidth = 9.86
$ws.Columns("I:J").ColumnWidth = 10.29
$ws.Columns("L").ColumnWidth = 17.86

Write-Output "done"
